# Corrections to perf numbers.
# Updates the "Atomix" worksheet (A1:D18 range) with corrected test names /
# performance figures, and moves the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atomix")

# Final (corrected) A:D data for rows 2-18.
$data = @(
    @(2,  "testTX9Mbps",  83.485799999999998, 31.994,             38.157499999999999),
    @(3,  "testTX6Mbps",  60.640500000000003, 25.4939,            28.790400000000002),
    @(4,  "testTX54Mbps", 221.26400000000001, 79.056799999999996, 97.892099999999999),
    @(5,  "testTX48Mbps", 230.12200000000001, 111.008,            165.69800000000001),
    @(6,  "testTX36Mbps", 234.87899999999999, 72.156400000000005, 88.995800000000003),
    @(7,  "testTX24Mbps", 197.06899999999999, 92.458799999999997, 148.13300000000001),
    @(8,  "testTX18Mbps", 141.738,             63.415999999999997, 79.016199999999998),
    @(9,  "testTX12Mbps", 112.996,             68.573800000000006, 96.390699999999995),
    @(10, "testRXCCA",    234.142,             140.99799999999999, 202.429),
    @(11, "testRX9Mbps",  119.5,               82.602599999999995, 114.246),
    @(12, "testRX6Mbps",  130.785,             92.693399999999997, 126.254),
    @(13, "testRX54Mbps", 50.476500000000001, 40.002400000000002, 39.424100000000003),
    @(14, "testRX48Mbps", 50.630499999999998, 38.575299999999999, 48.1678),
    @(15, "testRX36Mbps", 53.211399999999998, 40.762700000000002, 48.375900000000001),
    @(16, "testRX24Mbps", 67.674700000000001, 48.226300000000002, 61.790999999999997),
    @(17, "testRX18Mbps", 84.192999999999998, 63.304099999999998, 82.345399999999998),
    @(18, "testRX12Mbps", 99.204300000000003, 70.402500000000003, 79.364999999999995)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Move the active cell / selection on the Atomix sheet.
$ws.Activate()
$ws.Range("A22").Select()
